# 7.Matrix_machine.docx edit:
#   "shell is restricted except vi" paragraph gets extended with extra
#   instructions, and a brand-new paragraph ("go to vi ") is inserted
#   right after it.
#
# (The rest of the upstream diff is Word's automatic spell/grammar-check
#  <w:proofErr/> run-splitting and a namespace/latent-style cleanup that
#  a real Word session adds on save - it carries no textual content
#  change, so it is not reproduced here.)

$d = $word.ActiveDocument

# --- 1) Grow the "shell is restricted except vi" paragraph -----------------
# Locate it by its distinctive (pre-edit) text.
$p = $d.Paragraphs(15)
$r = $p.Range
[void]$r.MoveEnd(1, -1)          # exclude the paragraph mark from the range

# Append a throw-away trailing character ("X") together with the real new
# text. Placing the (collapsed) _GoBack bookmark immediately before that
# still-present dummy character -- rather than exactly at the then-last
# character of the run -- reliably anchors it at the run's end; deleting
# the dummy character afterwards leaves the bookmark sitting right after
# the final run, matching Word's own placement.
$r.Text = "shell is restricted except vieditor esc :!/bin/bassh  (restricyon will remove)X"

# --- 2) Re-add the "_GoBack" bookmark Word leaves at the last edit --------
$bmRange = $d.Range($r.End - 1, $r.End - 1)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the dummy trailing character now that the bookmark is anchored.
$p = $d.Paragraphs(15)
$delRange = $p.Range
[void]$delRange.MoveEnd(1, -1)
[void]$delRange.MoveStart(1, $delRange.End - $delRange.Start - 1)
$delRange.Text = ""

# --- 3) Insert the new "go to vi " paragraph right after it ---------------
$p = $d.Paragraphs(15)
[void]$p.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(16)
$newPara.Range.Text = "go to vi "
